# Updates cryptos list values (Price/Volume columns, plus the Filecoin /
# InjectiveProtocol row swap) to match the refreshed data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.491.28'
$ws.Range("E2").Value = '  +3.65%  '

$ws.Range("D3").Value = '3.498.33'
$ws.Range("E3").Value = '  +2.18%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''590.14'
$ws.Range("E5").Value = '  +3.16%  '

$ws.Range("D6").Value = '''169.50'
$ws.Range("E6").Value = '  +5.17%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.494.96'
$ws.Range("E8").Value = '  +2.01%  '

$ws.Range("D9").Value = '''0.588'
$ws.Range("E9").Value = '  +6.10%  '

$ws.Range("E10").Value = '  +0.51%  '

$ws.Range("E11").Value = '  +4.74%  '

$ws.Range("E12").Value = '  +3.17%  '

$ws.Range("D13").Value = '4.105.21'
$ws.Range("E13").Value = '  +2.31%  '

$ws.Range("E14").Value = '  -0.53%  '

$ws.Range("D15").Value = '''28.27'
$ws.Range("E15").Value = '  +4.63%  '

$ws.Range("E16").Value = '  +1.78%  '

$ws.Range("D17").Value = '66.516.94'
$ws.Range("E17").Value = '  +3.63%  '

$ws.Range("D18").Value = '3.514.17'
$ws.Range("E18").Value = '  +3.02%  '

$ws.Range("E19").Value = '  +3.92%  '

$ws.Range("E20").Value = '  +3.67%  '

$ws.Range("D21").Value = '''390.44'
$ws.Range("E21").Value = '  +3.59%  '

$ws.Range("D22").Value = '''7.97'
$ws.Range("E22").Value = '  +2.08%  '

$ws.Range("D23").Value = '''72.96'
$ws.Range("E23").Value = '  +2.37%  '

$ws.Range("D24").Value = '''0.998'
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("E25").Value = '  +3.28%  '

$ws.Range("E26").Value = '  +5.41%  '

$ws.Range("D27").Value = '''10.49'
$ws.Range("E27").Value = '  +10.15%  '

$ws.Range("E28").Value = '  +2.49%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("E30").Value = '  +5.50%  '

$ws.Range("D31").Value = '''1.48'
$ws.Range("E31").Value = '  +5.89%  '

$ws.Range("D32").Value = '''2.06'
$ws.Range("E32").Value = '  +2.60%  '

$ws.Range("D33").Value = '''23.58'
$ws.Range("E33").Value = '  +3.10%  '

$ws.Range("D34").Value = '''7.40'
$ws.Range("E34").Value = '  +4.30%  '

$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("E36").Value = '  +7.10%  '

$ws.Range("D37").Value = '''162.29'
$ws.Range("E37").Value = '  +1.64%  '

$ws.Range("D38").Value = '''0.882'
$ws.Range("E38").Value = '  +2.88%  '

$ws.Range("E39").Value = '  +4.83%  '

$ws.Range("D40").Value = '''6.89'
$ws.Range("E40").Value = '  +7.05%  '

$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''27.67'
$ws.Range("E41").Value = '  +6.51%  '

$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '''4.69'
$ws.Range("E42").Value = '  +6.04%  '

$ws.Range("E43").Value = '  +2.73%  '

$ws.Range("D44").Value = '''26.43'
$ws.Range("E44").Value = '  +2.40%  '

$ws.Range("D45").Value = '2.798.78'
$ws.Range("E45").Value = '  -0.21%  '

$ws.Range("E46").Value = '  +0.52%  '

$ws.Range("E47").Value = '  +1.95%  '

$ws.Range("E48").Value = '  +3.68%  '

$ws.Range("D49").Value = '''353.11'
$ws.Range("E49").Value = '  +5.42%  '

$ws.Range("E50").Value = '  +3.25%  '

$ws.Range("E51").Value = '  +12.33%  '
